$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row of raw/clean data for August 5th (row 67).
# Column A holds a date-formatted string ("2020-08-05"); writing it
# directly via .Value would make Excel auto-detect it as a real date
# and reformat the cell. Using a text formula and then collapsing it
# to a static value keeps it as plain text without changing the cell
# style.
$ws.Range("A67").Formula = "=""2020-08-05"""
$ws.Range("A67").Copy()
$ws.Range("A67").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B67").Value = 456100
$ws.Range("C67").Value = 499915
$ws.Range("D67").Value = 85845
$ws.Range("E67").Value = 49698
$ws.Range("F67").Value = 26.83
